$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.722.17"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.083.20"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.75"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.25"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.078.17"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.74"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "3.581.51"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "63.692.03"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "3.087.74"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.70"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.84"
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.18"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.84"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.24"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.15"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.40"
$ws.Range("E33").Value = "  -5.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.25"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "501.07"
$ws.Range("E36").Value = "  -4.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "3.266.19"
$ws.Range("E38").Value = "  +6.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0799"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.16"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.258"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("D47").Value = "0.0₃0540"
$ws.Range("E47").Value = "  +5.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.03"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.29"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  -13.94%  "
